$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 08:16:42"
$wsZhCn.Range("E5").Value = "2016-03-18 08:16:42"
$wsZhCn.Range("H2").Value = "2016-03-18 08:17:00"
$wsZhCn.Range("H5").Value = "2016-03-18 08:17:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 08:16:45"
$wsDeDe.Range("E5").Value = "2016-03-18 08:16:45"
$wsDeDe.Range("H2").Value = "2016-03-18 08:17:07"
$wsDeDe.Range("H5").Value = "2016-03-18 08:17:07"
